# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet holds one pitching-log row per game (rows 2-66), with column G
# labelled "K" (strikeouts). This recomputes/rewrites the K values for every
# row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") counts, keyed by worksheet row number.
$kValues = [ordered]@{
    2  = 1;  3  = 0;  4  = 1;  5  = 1;  6  = 2;  7  = 0;  8  = 2;  9  = 0;
    10 = 3;  11 = 0;  12 = 1;  13 = 3;  14 = 0;  15 = 2;  16 = 1;  17 = 1;
    18 = 2;  19 = 1;  20 = 0;  21 = 0;  22 = 1;  23 = 1;  24 = 1;  25 = 1;
    26 = 2;  27 = 2;  28 = 3;  29 = 0;  30 = 1;  31 = 0;  32 = 1;  33 = 3;
    34 = 3;  35 = 0;  36 = 2;  37 = 1;  38 = 3;  39 = 1;  40 = 0;  41 = 2;
    42 = 3;  43 = 2;  44 = 3;  45 = 1;  46 = 3;  47 = 1;  48 = 0;  49 = 3;
    50 = 2;  51 = 4;  52 = 1;  53 = 1;  54 = 1;  55 = 1;  56 = 2;  57 = 1;
    58 = 3;  59 = 3;  60 = 4;  61 = 2;  62 = 0;  63 = 2;  64 = 2;  65 = 2;
    66 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
